$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" column header in H1, reusing the same formatting
# (bold font, border, centered alignment) already used for the other
# header cells (e.g. G1) by copying its format over.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the corresponding data value in H2.
$ws.Range("H2").Value = 1
